# feat: Added province field
#
# Adds a new "PROVINCE" column (D) to the worksheet, populated with the
# same JSON value for every data row, and updates the selection/dimension
# accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("D1").Value = "PROVINCE"

# Same province value repeated for each of the 3 data rows
$province = '{"nome":"Luanda","id": 1}'
$ws.Range("D2").Value = $province
$ws.Range("D3").Value = $province
$ws.Range("D4").Value = $province

# Size the new column to fit its content (closest achievable width to the
# authored workbook's 22.7265625 character-width column D)
$ws.Columns.Item(4).ColumnWidth = 21.83

# Move/restore the active selection, as recorded in the authored workbook
$ws.Range("E9").Select()
